# Poland I Liga - atualizacao de bases (15-06-2024)
#
# This edit reshuffles a handful of match rows so that each row's fixture
# data (id, away team, score, odds, ...) lines up against the correct
# (fixed-in-place) home team for that row, and swaps a few duplicated
# fixture rows so the earlier row holds the earlier-id fixture.
#
# Concretely (worked out by diffing the underlying OOXML):
#   * Rows 3..10 rotate: row 3 receives what used to be row 10's fixture
#     data, and rows 4..10 each receive what used to be the previous row's
#     fixture data. Column A (running id) and column E (home team) stay
#     exactly where they are; everything else (B, F..AD) moves.
#   * Row pairs (75,76), (221,222), (272,273), (282,283) each have their
#     full fixture swapped between the two rows (column A stays put; B,
#     E, F..AD all swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row, $firstCol, $lastCol) {
    # Returns a 1-D array of the values from firstCol..lastCol (inclusive)
    # on the given row, read one cell at a time (safe/simple for a single row).
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += , $ws.Cells.Item($row, $c).Value
    }
    return , $vals
}

function Set-RowData($row, $firstCol, $lastCol, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$i]
        $i++
    }
}

# Column numbers: A=1 B=2 C=3 D=4 E=5 F=6 G=7 ... AD=30
$colB  = 2
$colF  = 6
$colAD = 30
$colE  = 5

# ---------------------------------------------------------------------
# 1) Rotate rows 3..10 (B column, and F..AD columns). A and E untouched.
# ---------------------------------------------------------------------
$rotRows = @(3, 4, 5, 6, 7, 8, 9, 10)

# Snapshot B (id) and F..AD (away team + stats/odds) for each row first,
# since we are about to overwrite them.
$snapB = @{}
$snapFAD = @{}
foreach ($r in $rotRows) {
    $snapB[$r] = $ws.Cells.Item($r, $colB).Value
    $snapFAD[$r] = Get-RowData -row $r -firstCol $colF -lastCol $colAD
}

# new[3] = old[10]; new[r] = old[r-1] for r = 4..10
$ws.Cells.Item(3, $colB).Value = $snapB[10]
Set-RowData -row 3 -firstCol $colF -lastCol $colAD -vals $snapFAD[10]

for ($r = 10; $r -ge 4; $r--) {
    $ws.Cells.Item($r, $colB).Value = $snapB[$r - 1]
    Set-RowData -row $r -firstCol $colF -lastCol $colAD -vals $snapFAD[$r - 1]
}

# ---------------------------------------------------------------------
# 2) Swap complete fixture rows (B, E, F..AD) for these pairs.
# ---------------------------------------------------------------------
$swapPairs = @(
    @(75, 76),
    @(221, 222),
    @(272, 273),
    @(282, 283)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $b1 = $ws.Cells.Item($r1, $colB).Value
    $b2 = $ws.Cells.Item($r2, $colB).Value

    $e1 = $ws.Cells.Item($r1, $colE).Value
    $e2 = $ws.Cells.Item($r2, $colE).Value

    $fad1 = Get-RowData -row $r1 -firstCol $colF -lastCol $colAD
    $fad2 = Get-RowData -row $r2 -firstCol $colF -lastCol $colAD

    $ws.Cells.Item($r1, $colB).Value = $b2
    $ws.Cells.Item($r2, $colB).Value = $b1

    $ws.Cells.Item($r1, $colE).Value = $e2
    $ws.Cells.Item($r2, $colE).Value = $e1

    Set-RowData -row $r1 -firstCol $colF -lastCol $colAD -vals $fad2
    Set-RowData -row $r2 -firstCol $colF -lastCol $colAD -vals $fad1
}
